$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.253.42'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '2.519.34'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'304.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").Value = "'96.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = "'0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("D10").Value = "'36.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.0806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = "'7.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = "'0.112"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").Value = '2.905.05'
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '2.515.80'
$ws.Range("E15").Value = '  -3.14%  '
$ws.Range("D16").Value = "'15.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.65%  '
$ws.Range("D17").Value = "'0.858"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("D18").Value = '42.298.06'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").Value = '0.0₃0973'
$ws.Range("D21").Value = "'6.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("D22").Value = "'71.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").Value = "'251.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("E24").Value = '  -2.22%  '
$ws.Range("D25").Value = "'2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.81%  '
$ws.Range("D26").Value = "'26.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.61%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = "'2.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.89%  '
$ws.Range("D29").Value = "'10.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").Value = "'37.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").Value = "'5.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").Value = "'154.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("D34").Value = "'0.0787"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.77%  '
$ws.Range("E35").Value = '  -5.33%  '
$ws.Range("E36").Value = '  -5.12%  '
$ws.Range("D37").Value = "'18.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("E38").Value = '  +1.24%  '
$ws.Range("D39").Value = "'0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").Value = "'24.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = "'3.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").Value = "'3.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").Value = "'2.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").Value = "'0.0299"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.45%  '
$ws.Range("D46").Value = '2.038.94'
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("D47").Value = "'84.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("E48").Value = '  -3.98%  '
$ws.Range("D49").Value = '2.767.12'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").Value = "'101.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.37%  '
